$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "33.931.75"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").Value = "1.783.52"
$ws.Range("E3").Value = "  -0.34%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.43"
$ws.Range("E5").Value = "  +2.18%  "
$ws.Range("E6").Value = "  -1.57%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.11"
$ws.Range("E8").Value = "  -1.54%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.293"
$ws.Range("E9").Value = "  +3.25%  "
$ws.Range("E10").Value = "  -4.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0940"
$ws.Range("E11").Value = "  +1.37%  "
$ws.Range("D12").Value = "2.039.96"
$ws.Range("E12").Value = "  -0.31%  "
$ws.Range("E13").Value = "  +2.10%  "
$ws.Range("D14").Value = "1.781.71"
$ws.Range("E14").Value = "  +0.83%  "
$ws.Range("D15").Value = "33.908.48"
$ws.Range("E15").Value = "  -0.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.615"
$ws.Range("E16").Value = "  -1.85%  "
$ws.Range("E17").Value = "  +0.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.52"
$ws.Range("E18").Value = "  -0.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.94"
$ws.Range("E19").Value = "  -0.97%  "
$ws.Range("D20").Value = "0.0₃0770"
$ws.Range("E20").Value = "  -1.49%  "
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.65"
$ws.Range("E22").Value = "  -1.26%  "
$ws.Range("E23").Value = "  -0.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.06"
$ws.Range("E24").Value = "  -1.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "161.83"
$ws.Range("E25").Value = "  +2.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.14"
$ws.Range("E26").Value = "  +1.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.14"
$ws.Range("E27").Value = "  -1.37%  "
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("E29").Value = "  +0.30%  "
$ws.Range("E30").Value = "  +2.64%  "
$ws.Range("E31").Value = "  -1.12%  "
$ws.Range("E32").Value = "  -1.23%  "
$ws.Range("E33").Value = "  +1.94%  "
$ws.Range("E34").Value = "  +1.13%  "
$ws.Range("D35").Value = "1.394.72"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.647"
$ws.Range("E36").Value = "  +1.08%  "
$ws.Range("E37").Value = "  -1.24%  "
$ws.Range("E38").Value = "  +1.31%  "
$ws.Range("E39").Value = "  +8.33%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "79.69"
$ws.Range("E40").Value = "  +0.16%  "
$ws.Range("E41").Value = "  +0.48%  "
$ws.Range("E42").Value = "  -0.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "13.62"
$ws.Range("E43").Value = "  +13.55%  "
$ws.Range("E44").Value = "  -1.53%  "
$ws.Range("D45").Value = "0.0₆0140"
$ws.Range("E45").Value = "  +9.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0510"
$ws.Range("E46").Value = "  +3.37%  "
$ws.Range("E47").Value = "  +2.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.92"
$ws.Range("E48").Value = "  +0.80%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "107.53"
$ws.Range("E49").Value = "  +0.18%  "
$ws.Range("D50").Value = "1.942.17"
$ws.Range("E50").Value = "  -0.47%  "
$ws.Range("E51").Value = "  +0.21%  "
